# Insert a new weekly price record at row 53 (Macroferia Regional de Talca - Mango),
# shifting all subsequent rows down by one (old row 53 -> 54, ..., old row 116 -> 117).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 53:116 down to 54:117, leaving row 53 empty (but carrying the date
# column's number format down with it, same as Excel's native row-insert).
$ws.Rows("53:53").Insert()

# Populate the newly inserted row 53 with the new data point.
$ws.Range("A53").Value = 5
$ws.Range("B53").Value = "Macroferia Regional de Talca"
$ws.Range("C53").Value = "Maule"
$ws.Range("D53").Value = 44629
$ws.Range("E53").Value = 7
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100108
$ws.Range("H53").Value = "Tropicales y subtropicales"
$ws.Range("I53").Value = 100108002
$ws.Range("J53").Value = "Mango"
$ws.Range("K53").Value = "Sin especificar"
$ws.Range("L53").Value = "Primera"
$ws.Range("M53").Value = 220
$ws.Range("N53").Value = 6500
$ws.Range("O53").Value = 7000
$ws.Range("P53").Value = 6773
$ws.Range("Q53").Value = "`$/bandeja 4 kilos"
$ws.Range("R53").Value = "Perú"
$ws.Range("S53").Value = 1693
$ws.Range("T53").Value = 4
